$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in supplier/contractor name strings (comma -> period) ---
$nameFixes = @(
  @{Cell="E35";  Value="FERNANDEZ MARIO H. GALLICET OSCAR M"},
  @{Cell="E78";  Value="FERNANDEZ MARIO H. GALLICET OSCAR M"},
  @{Cell="E148"; Value="FERNANDEZ MARIO H. GALLICET OSCAR M"},
  @{Cell="E74";  Value="PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"},
  @{Cell="F74";  Value="PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"},
  @{Cell="E79";  Value="IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"},
  @{Cell="F79";  Value="IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"},
  @{Cell="E111"; Value="IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"},
  @{Cell="F111"; Value="IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"},
  @{Cell="E149"; Value="IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"},
  @{Cell="F149"; Value="IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"},
  @{Cell="E83";  Value="MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"},
  @{Cell="E150"; Value="MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"}
)

foreach ($fix in $nameFixes) {
  $ws.Range($fix.Cell).Value = $fix.Value
}

# --- Reformat amounts in column H from "1.234,56" (es-AR) to "1234.56" (plain) ---
# These are stored as text (shared strings), not numbers, in the source workbook,
# so force text entry (NumberFormat "@") while assigning, then restore the
# cell's default/general style so no stray formatting is left behind.
$amountFixes = @(
  @{Cell="H2"; Value="918.00"},
  @{Cell="H3"; Value="236.00"},
  @{Cell="H4"; Value="209.00"},
  @{Cell="H5"; Value="5665.00"},
  @{Cell="H6"; Value="70250.00"},
  @{Cell="H7"; Value="747.98"},
  @{Cell="H8"; Value="1104.48"},
  @{Cell="H9"; Value="348.00"},
  @{Cell="H10"; Value="45.00"},
  @{Cell="H11"; Value="170.00"},
  @{Cell="H12"; Value="878.17"},
  @{Cell="H13"; Value="17710.00"},
  @{Cell="H14"; Value="80371.37"},
  @{Cell="H15"; Value="8354.25"},
  @{Cell="H16"; Value="13676.61"},
  @{Cell="H17"; Value="5921.14"},
  @{Cell="H18"; Value="607.31"},
  @{Cell="H19"; Value="3772.60"},
  @{Cell="H20"; Value="2576.24"},
  @{Cell="H21"; Value="883.47"},
  @{Cell="H22"; Value="5900.76"},
  @{Cell="H23"; Value="88.00"},
  @{Cell="H24"; Value="677.50"},
  @{Cell="H25"; Value="298.21"},
  @{Cell="H26"; Value="7108.15"},
  @{Cell="H27"; Value="6303.88"},
  @{Cell="H28"; Value="3783.50"},
  @{Cell="H29"; Value="1935.30"},
  @{Cell="H30"; Value="74.48"},
  @{Cell="H31"; Value="875.00"},
  @{Cell="H32"; Value="551.00"},
  @{Cell="H33"; Value="963.90"},
  @{Cell="H34"; Value="6158.00"},
  @{Cell="H35"; Value="1349.00"},
  @{Cell="H36"; Value="1976.49"},
  @{Cell="H37"; Value="2888.40"},
  @{Cell="H38"; Value="22930.00"},
  @{Cell="H39"; Value="5584.65"},
  @{Cell="H40"; Value="17593.78"},
  @{Cell="H41"; Value="50.00"},
  @{Cell="H42"; Value="3626.30"},
  @{Cell="H43"; Value="375.00"},
  @{Cell="H44"; Value="2860.00"},
  @{Cell="H45"; Value="13075.60"},
  @{Cell="H46"; Value="737.35"},
  @{Cell="H47"; Value="1183.37"},
  @{Cell="H48"; Value="74.00"},
  @{Cell="H49"; Value="929.00"},
  @{Cell="H50"; Value="2114.52"},
  @{Cell="H51"; Value="16482.39"},
  @{Cell="H52"; Value="1482.72"},
  @{Cell="H53"; Value="2.00"},
  @{Cell="H54"; Value="895.41"},
  @{Cell="H55"; Value="37485.18"},
  @{Cell="H56"; Value="5.66"},
  @{Cell="H57"; Value="81.50"},
  @{Cell="H58"; Value="126.00"},
  @{Cell="H59"; Value="1239.99"},
  @{Cell="H60"; Value="1387.18"},
  @{Cell="H61"; Value="285.00"},
  @{Cell="H62"; Value="194.55"},
  @{Cell="H63"; Value="174.00"},
  @{Cell="H64"; Value="91.82"},
  @{Cell="H65"; Value="9996.00"},
  @{Cell="H66"; Value="164.30"},
  @{Cell="H67"; Value="2903.50"},
  @{Cell="H68"; Value="699.00"},
  @{Cell="H69"; Value="832.00"},
  @{Cell="H70"; Value="68.79"},
  @{Cell="H71"; Value="3790.00"},
  @{Cell="H72"; Value="16800.00"},
  @{Cell="H73"; Value="1300.00"},
  @{Cell="H74"; Value="410.00"},
  @{Cell="H75"; Value="259.56"},
  @{Cell="H76"; Value="10062.00"},
  @{Cell="H77"; Value="2333.15"},
  @{Cell="H78"; Value="1049.00"},
  @{Cell="H79"; Value="1193.85"},
  @{Cell="H80"; Value="5820.60"},
  @{Cell="H81"; Value="1750.00"},
  @{Cell="H82"; Value="325.00"},
  @{Cell="H83"; Value="710.00"},
  @{Cell="H84"; Value="79.00"},
  @{Cell="H85"; Value="300.00"},
  @{Cell="H86"; Value="450.00"},
  @{Cell="H87"; Value="40.00"},
  @{Cell="H88"; Value="548.00"},
  @{Cell="H89"; Value="5.70"},
  @{Cell="H90"; Value="231.52"},
  @{Cell="H91"; Value="971.25"},
  @{Cell="H92"; Value="456.30"},
  @{Cell="H93"; Value="20000.00"},
  @{Cell="H94"; Value="774.00"},
  @{Cell="H95"; Value="15.96"},
  @{Cell="H96"; Value="2354.00"},
  @{Cell="H97"; Value="1984.20"},
  @{Cell="H98"; Value="31.20"},
  @{Cell="H99"; Value="6952.11"},
  @{Cell="H100"; Value="76.63"},
  @{Cell="H101"; Value="28.18"},
  @{Cell="H102"; Value="25.00"},
  @{Cell="H103"; Value="25.00"},
  @{Cell="H104"; Value="163.06"},
  @{Cell="H105"; Value="2600.00"},
  @{Cell="H106"; Value="450.00"},
  @{Cell="H107"; Value="30.00"},
  @{Cell="H108"; Value="7471.00"},
  @{Cell="H109"; Value="912.00"},
  @{Cell="H110"; Value="18196.00"},
  @{Cell="H111"; Value="803.80"},
  @{Cell="H112"; Value="186.00"},
  @{Cell="H113"; Value="296.00"},
  @{Cell="H114"; Value="600.00"},
  @{Cell="H115"; Value="1030.00"},
  @{Cell="H116"; Value="150.00"},
  @{Cell="H117"; Value="480.00"},
  @{Cell="H118"; Value="100.00"},
  @{Cell="H119"; Value="440.00"},
  @{Cell="H120"; Value="447.60"},
  @{Cell="H121"; Value="1228.49"},
  @{Cell="H122"; Value="750.00"},
  @{Cell="H123"; Value="394.00"},
  @{Cell="H124"; Value="6.50"},
  @{Cell="H125"; Value="1766.41"},
  @{Cell="H126"; Value="450.00"},
  @{Cell="H127"; Value="240.00"},
  @{Cell="H128"; Value="500.00"},
  @{Cell="H129"; Value="150.00"},
  @{Cell="H130"; Value="1937.94"},
  @{Cell="H131"; Value="252.00"},
  @{Cell="H132"; Value="700.00"},
  @{Cell="H133"; Value="2781.75"},
  @{Cell="H134"; Value="380.00"},
  @{Cell="H135"; Value="550.00"},
  @{Cell="H136"; Value="8970.00"},
  @{Cell="H137"; Value="150.00"},
  @{Cell="H138"; Value="9824.76"},
  @{Cell="H139"; Value="1200.00"},
  @{Cell="H140"; Value="150.00"},
  @{Cell="H141"; Value="800.00"},
  @{Cell="H142"; Value="2500.00"},
  @{Cell="H143"; Value="500.00"},
  @{Cell="H144"; Value="350.00"},
  @{Cell="H145"; Value="120.00"},
  @{Cell="H146"; Value="3663.00"},
  @{Cell="H147"; Value="367.40"},
  @{Cell="H148"; Value="146.00"},
  @{Cell="H149"; Value="280.00"},
  @{Cell="H150"; Value="15.00"},
  @{Cell="H151"; Value="1895.00"},
  @{Cell="H152"; Value="40.00"},
  @{Cell="H153"; Value="2420.00"},
  @{Cell="H154"; Value="90.00"},
  @{Cell="H155"; Value="927.00"},
  @{Cell="H156"; Value="30.25"},
  @{Cell="H157"; Value="579.02"},
  @{Cell="H158"; Value="3125.00"},
  @{Cell="H159"; Value="774.00"},
  @{Cell="H160"; Value="1255.00"},
  @{Cell="H161"; Value="1670.01"},
  @{Cell="H162"; Value="294.06"},
  @{Cell="H163"; Value="337.50"},
  @{Cell="H164"; Value="300.00"},
  @{Cell="H165"; Value="1200.00"},
  @{Cell="H166"; Value="280.40"},
  @{Cell="H167"; Value="90.30"},
  @{Cell="H168"; Value="38.00"},
  @{Cell="H169"; Value="2795.05"},
  @{Cell="H170"; Value="670.00"},
  @{Cell="H171"; Value="390.00"},
  @{Cell="H172"; Value="2767.27"},
  @{Cell="H173"; Value="1240.82"},
  @{Cell="H174"; Value="4160.00"},
  @{Cell="H175"; Value="650.00"},
  @{Cell="H176"; Value="486232.47"},
  @{Cell="H177"; Value="121500.00"},
  @{Cell="H178"; Value="1019.30"},
  @{Cell="H179"; Value="5400.00"},
  @{Cell="H180"; Value="245.10"},
  @{Cell="H181"; Value="318.00"},
  @{Cell="H182"; Value="1400.00"}
)

foreach ($fix in $amountFixes) {
  $rng = $ws.Range($fix.Cell)
  $rng.NumberFormat = "@"
  $rng.Value = $fix.Value
  $rng.Style = "Normal"
}
